# Add task to run npm on start
# Update the experiment3 results sheet with refreshed cross-validation metrics.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("experiment3")

# Row 2: SVM
$ws.Range("B2").Value = 939.0947283049472
$ws.Range("C2").Value = 16.94216969952914
$ws.Range("D2").Value = 945.8715961847589
$ws.Range("E2").Value = 932.3178604251356

# Row 3: Linear Regression
$ws.Range("B3").Value = 1024.755594317984
$ws.Range("C3").Value = 5.720478823749245
$ws.Range("D3").Value = 1027.043785847484
$ws.Range("E3").Value = 1022.467402788484

# Row 4: KNN
$ws.Range("B4").Value = 1037.260178426602
$ws.Range("C4").Value = 7.640943522047623
$ws.Range("D4").Value = 1040.316555835421
$ws.Range("E4").Value = 1034.203801017783

# Row 5: MLPRegressor
$ws.Range("B5").Value = 1040.745979882588
$ws.Range("C5").Value = 6.515214841509487
$ws.Range("D5").Value = 1043.352065819191
$ws.Range("E5").Value = 1038.139893945984

# Row 6: GradientBoostingRegressor
$ws.Range("B6").Value = 1040.775288343738
$ws.Range("C6").Value = 3.880730399725143
$ws.Range("D6").Value = 1042.327580503628
$ws.Range("E6").Value = 1039.222996183848

# Row 7: now XGBRegressor (was RandomForestRegressor)
$ws.Range("A7").Value = "XGBRegressor"
$ws.Range("B7").Value = 1162.888669228534
$ws.Range("C7").Value = 15.68002153275733
$ws.Range("D7").Value = 1169.160677841637
$ws.Range("E7").Value = 1156.616660615431

# Row 8: now RandomForestRegressor (was XGBRegressor)
$ws.Range("A8").Value = "RandomForestRegressor"
$ws.Range("B8").Value = 1163.903413905478
$ws.Range("C8").Value = 18.58558246035076
$ws.Range("D8").Value = 1171.337646889619
$ws.Range("E8").Value = 1156.469180921338

# Row 9: RandomForestClassifier
$ws.Range("B9").Value = 1310.18199513382
$ws.Range("C9").Value = 43.25342831355597
$ws.Range("D9").Value = 1327.483366459242
$ws.Range("E9").Value = 1292.880623808398

# Row 10: DecisionTreeClassifier
$ws.Range("B10").Value = 1343.603892944039
$ws.Range("C10").Value = 33.86016428579097
$ws.Range("D10").Value = 1357.147958658355
$ws.Range("E10").Value = 1330.059827229722

# Row 11: Decision Tree Regressor
$ws.Range("B11").Value = 1363.088686131387
$ws.Range("C11").Value = 35.40871296532122
$ws.Range("D11").Value = 1377.252171317515
$ws.Range("E11").Value = 1348.925200945258

# Row 12: MLPClassifier
$ws.Range("B12").Value = 1424.411192214112
$ws.Range("C12").Value = 29.6972603791136
$ws.Range("D12").Value = 1436.290096365757
$ws.Range("E12").Value = 1412.532288062466
